# Atualizacao de bases das ligas, do dia: 14-05-2024 as 01:09
#
# 1) A handful of existing match rows had their B..AB data (id, teams,
#    odds, results, ...) corrected - each correction turns out to be an
#    exact swap of the B..AB payload between two neighbouring rows (the
#    sequential index in column A stays put).
# 2) Four brand-new matches are appended as rows 262-265.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold data (everything except the sequential index in A).
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

# Row pairs whose B..AB content needs to be swapped.
$pairs = @(
    @(91,92),
    @(95,96),
    @(192,193),
    @(200,201),
    @(231,232),
    @(237,238),
    @(249,250)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

# New matches to append after the current last row (261).
$newRows = @(
    @{ Row=262; A=260; B=8162888; C="Costa Rica Primera Division"; D=45420.75;          E="AD Guanacasteca";  F="Municipal Perez Zeledon"; G=4; H=0; I="H"; J=1.7;   K=3.5;  L=4;     M=1.666; N=3.8; O=4;    P=-0.75; Q=1.925; R=1.875; S=2.75; T=1.85; U=1.95; V=0.6659999999999999; W=-1; X=-1;  Y=0.925; Z=-1;    AA=0.8500000000000001; AB=-1 },
    @{ Row=263; A=261; B=8162889; C="Costa Rica Primera Division"; D=45420.875;          E="Cartagines";       F="Sporting San Jose";       G=1; H=1; I="D"; J=1.8;   K=3.4;  L=3.75;  M=2;     N=3.5; O=3.1;  P=-0.25; Q=1.8;   R=2;     S=2.75; T=1.8;  U=2;    V=-1;                 W=2.5; X=-1;  Y=-0.5;  Z=0.5;   AA=-1;                 AB=1  },
    @{ Row=264; A=262; B=8162887; C="Costa Rica Primera Division"; D=45420.91666666666;  E="Santos de Gupiles"; F="Alajuelense";             G=0; H=3; I="A"; J=5.2;   K=3.75; L=1.533; M=5.75;  N=4;   O=1.45; P=1;     Q=2.025; R=1.775; S=2.75; T=2;    U=1.8;  V=-1;                 W=-1;  X=0.45; Y=-1;    Z=0.7749999999999999;  AA=0.5;                AB=-0.5 },
    @{ Row=265; A=263; B=8162890; C="Costa Rica Primera Division"; D=45420.97916666666;  E="Herediano";        F="Puntarenas";              G=3; H=0; I="H"; J=1.3;   K=4.75; L=7.5;   M=1.181; N=6.5; O=9.5;  P=-2;    Q=2.025; R=1.775; S=3.25; T=2;    U=1.8;  V=0.181;              W=-1;  X=-1;  Y=1.025; Z=-1;    AA=-0.5;               AB=0.4 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    $ws.Range("A$r").Value2 = $nr.A
    $ws.Range("B$r").Value2 = $nr.B
    $ws.Range("C$r").Value2 = $nr.C
    $ws.Range("D$r").Value2 = $nr.D
    $ws.Range("E$r").Value2 = $nr.E
    $ws.Range("F$r").Value2 = $nr.F
    $ws.Range("G$r").Value2 = $nr.G
    $ws.Range("H$r").Value2 = $nr.H
    $ws.Range("I$r").Value2 = $nr.I
    $ws.Range("J$r").Value2 = $nr.J
    $ws.Range("K$r").Value2 = $nr.K
    $ws.Range("L$r").Value2 = $nr.L
    $ws.Range("M$r").Value2 = $nr.M
    $ws.Range("N$r").Value2 = $nr.N
    $ws.Range("O$r").Value2 = $nr.O
    $ws.Range("P$r").Value2 = $nr.P
    $ws.Range("Q$r").Value2 = $nr.Q
    $ws.Range("R$r").Value2 = $nr.R
    $ws.Range("S$r").Value2 = $nr.S
    $ws.Range("T$r").Value2 = $nr.T
    $ws.Range("U$r").Value2 = $nr.U
    $ws.Range("V$r").Value2 = $nr.V
    $ws.Range("W$r").Value2 = $nr.W
    $ws.Range("X$r").Value2 = $nr.X
    $ws.Range("Y$r").Value2 = $nr.Y
    $ws.Range("Z$r").Value2 = $nr.Z
    $ws.Range("AA$r").Value2 = $nr.AA
    $ws.Range("AB$r").Value2 = $nr.AB

    # Match the look of the existing rows: bold/centered/bordered index
    # cell in column A, and the date-time number format in column D.
    $ws.Range("A261").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("D261").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0
